$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 415.4
$ws.Range("I53").Value = 91.666664
$ws.Range("K53").Value = 91.666664
$ws.Range("M53").Value = 545.333336
$ws.Range("H80").Value = 395.5
$ws.Range("J80").Value = 507
$ws.Range("L80").Value = 1521
$ws.Range("N80").Value = -3517
$ws.Range("H83").Value = 395.5
$ws.Range("J83").Value = 507
$ws.Range("L83").Value = 4563
$ws.Range("N83").Value = -14547
$ws.Range("H86").Value = 4481
$ws.Range("I86").Value = 4444
$ws.Range("J86").Value = 4499.5
$ws.Range("K86").Value = 4444
$ws.Range("L86").Value = 4499.5
$ws.Range("M86").Value = -3321
$ws.Range("N86").Value = -6745.5
$ws.Range("H88").Value = 1331.6666
$ws.Range("I88").Value = 1497
$ws.Range("J88").Value = 1249
$ws.Range("K88").Value = 1497
$ws.Range("L88").Value = 1249
$ws.Range("M88").Value = -1091
$ws.Range("N88").Value = -2061
$ws.Range("H89").Value = 4481
$ws.Range("I89").Value = 4444
$ws.Range("J89").Value = 4499.5
$ws.Range("K89").Value = 22220
$ws.Range("L89").Value = 22497.5
$ws.Range("M89").Value = -16604
$ws.Range("N89").Value = -33729.5
$ws.Range("H91").Value = 1331.6666
$ws.Range("I91").Value = 1497
$ws.Range("J91").Value = 1249
$ws.Range("K91").Value = 1497
$ws.Range("L91").Value = 1249
$ws.Range("M91").Value = -93
$ws.Range("N91").Value = -4057
$ws.Range("H92").Value = 581
$ws.Range("I92").Value = 496.5
$ws.Range("J92").Value = 665.5
$ws.Range("K92").Value = 496.5
$ws.Range("L92").Value = 665.5
$ws.Range("M92").Value = 751.5
$ws.Range("N92").Value = -3161.5
$ws.Range("H96").Value = 1005.8
$ws.Range("I96").Value = 876.6667
$ws.Range("J96").Value = 1199.5
$ws.Range("K96").Value = 2630.0001
$ws.Range("L96").Value = 3598.5
$ws.Range("M96").Value = -1257.0001
$ws.Range("N96").Value = -6344.5
$ws.Range("H97").Value = 399
$ws.Range("J97").Value = 399
$ws.Range("L97").Value = 1197
$ws.Range("N97").Value = -2189
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 1525.8572
$ws.Range("I132").Value = 1577.8462
$ws.Range("J132").Value = 850
$ws.Range("K132").Value = 4733.5386
$ws.Range("L132").Value = 2550
$ws.Range("M132").Value = -2203.5386
$ws.Range("N132").Value = -7610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 35997
$ws.Range("J23").Value = 35997
$ws.Range("L23").Value = 35997
$ws.Range("N23").Value = -36515
$ws.Range("H32").Value = 4123.4
$ws.Range("I32").Value = 4146.952
$ws.Range("J32").Value = 3999.75
$ws.Range("K32").Value = 4146.952
$ws.Range("L32").Value = 3999.75
$ws.Range("M32").Value = -3859.952
$ws.Range("N32").Value = -4573.75
$ws.Range("H37").Value = 35998.08
$ws.Range("J37").Value = 35998.08
$ws.Range("L37").Value = 35998.08
$ws.Range("N37").Value = -36544.08
$ws.Range("H74").Value = 1319.8889
$ws.Range("I74").Value = 1411.2858
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1411.2858
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -537.2858000000001
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 1319.8889
$ws.Range("I77").Value = 1411.2858
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 7056.429
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -2688.429
$ws.Range("N77").Value = -13736
$ws.Range("H122").Value = 924.7368
$ws.Range("I122").Value = 892.8333
$ws.Range("K122").Value = 2678.4999
$ws.Range("M122").Value = -228.4998999999998
$ws.Range("H132").Value = 1644.0264
$ws.Range("I132").Value = 1370.4546
$ws.Range("K132").Value = 4111.3638
$ws.Range("M132").Value = -1581.3638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3777.0557
$ws.Range("I134").Value = 3330.0667
$ws.Range("K134").Value = 9990.2001
$ws.Range("M134").Value = -7455.2001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1819.25
$ws.Range("I31").Value = 1864
$ws.Range("K31").Value = 1864
$ws.Range("M31").Value = -1569
$ws.Range("H34").Value = 1819.25
$ws.Range("I34").Value = 1864
$ws.Range("K34").Value = 1864
$ws.Range("M34").Value = -1662

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 334.8889
$ws.Range("I12").Value = 256
$ws.Range("K12").Value = 768
$ws.Range("M12").Value = -595
$ws.Range("H129").Value = 1677.8
$ws.Range("I129").Value = 1452.8334
$ws.Range("J129").Value = 2015.25
$ws.Range("K129").Value = 4358.5002
$ws.Range("L129").Value = 6045.75
$ws.Range("M129").Value = 641.4997999999996
$ws.Range("N129").Value = -16045.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4359
$ws.Range("I80").Value = 2983
$ws.Range("J80").Value = 5735
$ws.Range("K80").Value = 2983
$ws.Range("L80").Value = 5735
$ws.Range("M80").Value = -1985
$ws.Range("N80").Value = -7731
$ws.Range("H83").Value = 4359
$ws.Range("I83").Value = 2983
$ws.Range("J83").Value = 5735
$ws.Range("K83").Value = 14915
$ws.Range("L83").Value = 28675
$ws.Range("M83").Value = -9923
$ws.Range("N83").Value = -38659
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3616.1667
$ws.Range("I132").Value = 2924.75
$ws.Range("K132").Value = 8774.25
$ws.Range("M132").Value = -6244.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4799.4
$ws.Range("J46").Value = 8998.5
$ws.Range("L46").Value = 8998.5
$ws.Range("N46").Value = -9374.5
$ws.Range("H55").Value = 212.5
$ws.Range("J55").Value = 266.33334
$ws.Range("L55").Value = 266.33334
$ws.Range("N55").Value = -612.33334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H81").Value = 5100.857
$ws.Range("I81").Value = 2450.6667
$ws.Range("K81").Value = 4901.3334
$ws.Range("M81").Value = -3840.3334
$ws.Range("H84").Value = 5100.857
$ws.Range("I84").Value = 2450.6667
$ws.Range("K84").Value = 24506.667
$ws.Range("M84").Value = -19202.667
$ws.Range("H100").Value = 689.5
$ws.Range("I100").Value = 689
$ws.Range("K100").Value = 1378
$ws.Range("M100").Value = -837
$ws.Range("H132").Value = 2028.5151
$ws.Range("I132").Value = 1366.1818
$ws.Range("K132").Value = 4098.5454
$ws.Range("M132").Value = -1568.5454

